# feat: add 2022-Q1 data
#
# 1. The existing "总计" (summary) worksheet becomes the new "2022-Q1"
#    fund-holdings detail sheet (same slot/position as before, right
#    after "2021-Q4").
# 2. A brand-new "总计" worksheet is appended at the end, holding the
#    summary table with a new top row for "2022-Q1", on top of the
#    previous "2021-Q4"/"2021-Q3" rows.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Helper to write a value as literal text (avoids Excel's automatic
# "looks like a number" coercion) while keeping the cell's style at the
# workbook default ("Normal" / no explicit xf), matching the source data.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet as the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Reuse the header / first-column formatting (bold, centered, bordered
# style) from the "2021-Q4" sheet so the new sheet matches the look of
# the other quarterly detail sheets.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q4.Range("A2").Copy()
$q1.Range("A2:A4").PasteSpecial($xlPasteFormats)

$q1.Application.CutCopyMode = $false

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 004008
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "004008"
Set-TextValue $q1.Range("C2") "中融鑫思路灵活配置混合A"
Set-TextValue $q1.Range("D2") "3.91"
Set-TextValue $q1.Range("E2") "35.78"
Set-TextValue $q1.Range("F2") "1.54"
Set-TextValue $q1.Range("G2") "0.0602"
$q1.Range("H2").Value = 5

# Row 3 - 004009
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "004009"
Set-TextValue $q1.Range("C3") "中融鑫思路灵活配置混合C"
Set-TextValue $q1.Range("D3") "2.36"
Set-TextValue $q1.Range("E3") "35.78"
Set-TextValue $q1.Range("F3") "1.54"
Set-TextValue $q1.Range("G3") "0.0363"
$q1.Range("H3").Value = 5

# Row 4 - 501032
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "501032"
Set-TextValue $q1.Range("C4") "财通福盛多策略混合（LOF）"
Set-TextValue $q1.Range("D4") "0.74"
Set-TextValue $q1.Range("E4") "92.47"
Set-TextValue $q1.Range("F4") "4.11"
Set-TextValue $q1.Range("G4") "0.0304"
$q1.Range("H4").Value = 8

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" (summary) sheet at the end
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Reuse the header / first-column formatting from the "2021-Q4" sheet.
$q4.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)

$q4.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial($xlPasteFormats)

$total.Application.CutCopyMode = $false

# Header row
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Row 2 - 2022-Q1 (new)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.13

# Row 3 - 2021-Q4 (previously row 2)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 12
$total.Range("D3").Value = 0.96

# Row 4 - 2021-Q3 (previously row 3)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0
